# MPC boilers linearized: update the offset used in column C's formula
# (C<row> = B<row> + <offset>) for rows 103-130 of the hot water
# consumption profile sheet, and update the view's selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => new additive offset applied to column B to compute column C
$offsets = @{
    103 = 10
    104 = 10
    105 = 10
    106 = 10
    107 = 15
    108 = 15
    109 = 15
    110 = 15
    111 = 20
    112 = 20
    113 = 20
    114 = 20
    115 = 3
    116 = 3
    117 = 3
    118 = 3
    119 = 3
    120 = 3
    121 = 20
    122 = 20
    123 = 20
    124 = 20
    125 = -3
    126 = -3
    127 = -3
    128 = -3
    129 = -3
    130 = -3
}

foreach ($row in $offsets.Keys) {
    $offset = $offsets[$row]
    if ($offset -ge 0) {
        $formula = "=B$row+$offset"
    } else {
        $formula = "=B$row$offset"
    }
    $ws.Range("C$row").Formula = $formula
}

# Match the updated view state: scrolled down with C131 as the active cell
$ws.Range("C131").Select()
